# langtjern parameter updates from GIS shp
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LakePerimeter (row 3) updated from GIS shapefile
$ws.Range("B3").Value = 4221

# LakeArea (row 5) updated from GIS shapefile
$ws.Range("B5").Value = 224263

# Update the active selection to reflect where the editor left off
$ws.Range("E4").Select()
